# Refresh the cryptos price/volume snapshot (GitHub Actions scheduled update).
# Price cells (column D) that look like plain numbers ("215.02", "1.003", ...)
# are forced to Text format before assignment so Excel's COM layer doesn't
# silently coerce them to numeric values (which would drop the trailing
# zeros / change the stored type) - matching the original inline-string cells.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '25.815.84'
$ws.Range('E2').Value = '  -1.35%  '
$ws.Range('D3').Value = '1.634.10'
$ws.Range('E3').Value = '  -1.45%  '
$ws.Range('E4').Value = '  -0.25%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '215.02'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.26%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.5021'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -2.35%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.003'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.27%  '
$ws.Range('E8').Value = '  -0.37%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06386'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.72%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '19.59'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -1.80%  '
$ws.Range('E11').Value = '  -1.66%  '
$ws.Range('D12').Value = '1.639.51'
$ws.Range('E12').Value = '  -1.13%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.231'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -1.49%  '
$ws.Range('D14').Value = '1.858.88'
$ws.Range('E14').Value = '  -1.44%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.5454'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -1.68%  '
$ws.Range('D16').Value = '0.0₅7917'
$ws.Range('E16').Value = '  -1.73%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '63.52'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -1.10%  '
$ws.Range('D18').Value = '25.832.43'
$ws.Range('E18').Value = '  -1.47%  '
$ws.Range('E19').Value = '  -0.22%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '202.90'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -4.09%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '4.303'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -2.77%  '
$ws.Range('E22').Value = '  -1.13%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.965'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.18%  '
$ws.Range('E24').Value = '  -0.23%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '1.933'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +10.17%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '141.16'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E27').Value = '  -1.77%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '15.70'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -0.56%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '6.686'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -4.10%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.240'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -1.10%  '
$ws.Range('E31').Value = '  -4.85%  '
$ws.Range('E32').Value = '  -2.77%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.188'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.88%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.531'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -2.38%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.349'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -0.97%  '
$ws.Range('D36').Value = '1.171.90'
$ws.Range('E36').Value = '  -0.09%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.622'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -5.06%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.8909'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -4.24%  '
$ws.Range('E39').Value = '  -2.21%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.01555'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -2.41%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.565'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +0.08%  '
$ws.Range('E42').Value = '  -0.27%  '
$ws.Range('E43').Value = '  -0.48%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '99.34'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -1.19%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.8013'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -5.16%  '
$ws.Range('D46').Value = '1.771.45'
$ws.Range('E46').Value = '  -1.36%  '
$ws.Range('D47').Value = '0.0₈111'
$ws.Range('E47').Value = '  -2.85%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.4511'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.57%  '
$ws.Range('E49').Value = '  +0.19%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '54.75'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -2.07%  '
$ws.Range('E51').Value = '  -0.51%  '
